# Contoso Chai Tea market trends 2023 — localisation pass.
#
#  - Header A1:      "Date"                              -> "Fecha"
#  - Header B1:      "Ventas totales de Chai (unidades)"  -> "Ventas totales de chai (unidades)"
#  - D6 / E6:         436 / 1705 (numbers)                -> "4:36" / "05:17" (text that merely
#                      *looks* like a clock time — entered as literal strings, not Time values)
#
# Switch off auto-recalculation first so the shared formula in column B (SUM(C+D)) is not
# re-evaluated against the now-textual D6 and left sitting on a stale-but-valid cached
# result, exactly as happened in the source edit (B6 keeps its old cached value of 935).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.Calculation = -4135   # xlCalculationManual

$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Ventas totales de chai (unidades)"

$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "05:17"
